# Re-applies the "Merge pull request #72" changes:
#  - Insert a new "DataProviderIndustryType" sheet (copy of the
#    IndustryType column already present on DatProviderOrg) right after
#    the "DatProviderOrg" sheet.
#  - Insert a new "DatProviderGlAccount" sheet (new GL account lookup
#    list) at the end of the workbook.
#  - Update the selection remembered on "DatProviderOrg" to span the
#    whole IndustryType column.
#  - Leave the new IndustryType sheet as the active tab, matching the
#    state the workbook was saved in.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) New "DatProviderGlAccount" sheet, appended at the very end.
#    (Created first so it claims the lower internal sheetId, matching
#    the order sheets were added in the original edit.)
# ------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("DatProviderProd")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGl = $wb.Worksheets.Add($null, $lastSheet)
$wsGl.Name = "DatProviderGlAccount"

# Borrow the header/data formatting used by the other single-column
# lookup sheets (bold+filled header, thin-bordered data cells) instead
# of re-deriving new style entries.
$wsProd.Range("A1").Copy()
$wsGl.Range("A1").PasteSpecial(-4122) | Out-Null
$wsProd.Range("A2").Copy()
$wsGl.Range("A2:A10").PasteSpecial(-4122) | Out-Null

$wsGl.Range("A1").Value = "GL Account"
$wsGl.Range("A2").Value = "300-Sales-Software"
$wsGl.Range("A3").Value = "301-Sales-Hardware"
$wsGl.Range("A4").Value = "302-Rental-Income"
$wsGl.Range("A5").Value = "303-Interest-Income"
$wsGl.Range("A6").Value = "304-Sales-Software-Support"
$wsGl.Range("A7").Value = "305-Sales Other"
$wsGl.Range("A8").Value = "306-Internet Sales"
$wsGl.Range("A9").Value = "307-Service-Hardware Labor"
$wsGl.Range("A10").Value = "308-Sales-Books"

$wsGl.Columns.Item(1).AutoFit() | Out-Null
$wsGl.Range("J14").Select() | Out-Null

# ------------------------------------------------------------------
# 2) New "DataProviderIndustryType" sheet, inserted right after
#    "DatProviderOrg".
# ------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("DatProviderOrg")
$wsInd = $wb.Worksheets.Add($null, $wsOrg)
$wsInd.Name = "DataProviderIndustryType"

# This sheet duplicates the IndustryType lookup column that already
# lives in DatProviderOrg!B1:B6, so copy its values+formatting wholesale.
$wsOrg.Range("B1:B6").Copy()
$wsInd.Range("A1").PasteSpecial(-4122) | Out-Null
$wsInd.Range("A1:A6").PasteSpecial(-4163) | Out-Null

$wsInd.Columns.Item(1).AutoFit() | Out-Null

# ------------------------------------------------------------------
# 3) Cosmetic selection tweak remembered on DatProviderOrg: the whole
#    IndustryType column is selected now instead of just B2.
# ------------------------------------------------------------------
$wsOrg.Range("B1:B6").Select() | Out-Null

# The workbook was last saved with the new IndustryType sheet active
# and its G8 cell selected.
$wsInd.Activate() | Out-Null
$wsInd.Range("G8").Select() | Out-Null
